# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows at the top of the "Palta" data block
# (row 886), pushing the existing rows down by 3 (old row 886 -> 889,
# ..., old row 936 -> 939) and fill the 3 new rows with the latest
# week's prices (Especial / Primera / Segunda, Perú, $/bandeja 10 kilos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 886:936 down to 889:939 by inserting 3 blank rows.
$ws.Range("886:888").Insert()

# Row 886 - Especial
$ws.Cells.Item(886,1).Value = 8
$ws.Cells.Item(886,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(886,3).Value = "Coquimbo"
$ws.Cells.Item(886,4).Value = 44783
$ws.Cells.Item(886,5).Value = 4
$ws.Cells.Item(886,6).Value = "Fruta"
$ws.Cells.Item(886,7).Value = 100106
$ws.Cells.Item(886,8).Value = "Oleaginosos"
$ws.Cells.Item(886,9).Value = 100106002
$ws.Cells.Item(886,10).Value = "Palta"
$ws.Cells.Item(886,11).Value = "Hass"
$ws.Cells.Item(886,12).Value = "Especial"
$ws.Cells.Item(886,13).Value = 480
$ws.Cells.Item(886,14).Value = 28000
$ws.Cells.Item(886,15).Value = 29000
$ws.Cells.Item(886,16).Value = 28500
$ws.Cells.Item(886,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(886,18).Value = "Perú"
$ws.Cells.Item(886,19).Value = 2850
$ws.Cells.Item(886,20).Value = 10

# Row 887 - Primera
$ws.Cells.Item(887,1).Value = 8
$ws.Cells.Item(887,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(887,3).Value = "Coquimbo"
$ws.Cells.Item(887,4).Value = 44783
$ws.Cells.Item(887,5).Value = 4
$ws.Cells.Item(887,6).Value = "Fruta"
$ws.Cells.Item(887,7).Value = 100106
$ws.Cells.Item(887,8).Value = "Oleaginosos"
$ws.Cells.Item(887,9).Value = 100106002
$ws.Cells.Item(887,10).Value = "Palta"
$ws.Cells.Item(887,11).Value = "Hass"
$ws.Cells.Item(887,12).Value = "Primera"
$ws.Cells.Item(887,13).Value = 400
$ws.Cells.Item(887,14).Value = 25000
$ws.Cells.Item(887,15).Value = 26000
$ws.Cells.Item(887,16).Value = 25500
$ws.Cells.Item(887,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(887,18).Value = "Perú"
$ws.Cells.Item(887,19).Value = 2550
$ws.Cells.Item(887,20).Value = 10

# Row 888 - Segunda
$ws.Cells.Item(888,1).Value = 8
$ws.Cells.Item(888,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(888,3).Value = "Coquimbo"
$ws.Cells.Item(888,4).Value = 44783
$ws.Cells.Item(888,5).Value = 4
$ws.Cells.Item(888,6).Value = "Fruta"
$ws.Cells.Item(888,7).Value = 100106
$ws.Cells.Item(888,8).Value = "Oleaginosos"
$ws.Cells.Item(888,9).Value = 100106002
$ws.Cells.Item(888,10).Value = "Palta"
$ws.Cells.Item(888,11).Value = "Hass"
$ws.Cells.Item(888,12).Value = "Segunda"
$ws.Cells.Item(888,13).Value = 360
$ws.Cells.Item(888,14).Value = 23000
$ws.Cells.Item(888,15).Value = 24000
$ws.Cells.Item(888,16).Value = 23500
$ws.Cells.Item(888,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(888,18).Value = "Perú"
$ws.Cells.Item(888,19).Value = 2350
$ws.Cells.Item(888,20).Value = 10
